$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.567.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.00%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.660.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.32%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  +0.11%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'601.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.14%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'156.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.19%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  +0.09%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.622"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +5.57%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Value = "'  +0.06%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = "'  -1.02%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'5.86"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -2.24%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Value = "'  -0.26%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Value = "'  -2.49%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("E14").Value = "'  -5.11%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'3.136.89"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.16%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'65.407.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.04%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'2.661.38"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.45%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'12.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.25%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'4.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.21%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = "'  +2.02%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'351.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.10%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "'  -0.14%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'69.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.13%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = "'  +2.76%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'9.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.21%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "'  -4.27%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").Value = "'  -2.25%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'  -3.57%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'8.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.82%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").Value = "'  -0.20%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'529.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.01%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").Value = "'  -3.17%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'1.73"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -2.99%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'6.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.51%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'5.48"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.55%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "'  -2.47%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'20.34"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.79%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = "'  +0.17%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'159.28"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.05%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").Value = "'  -3.92%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E42").Value = "'  +1.76%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'164.45"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -3.15%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'4.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.96%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("E45").Value = "'  -0.36%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'2.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -2.32%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "'  -0.63%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").Value = "'  -2.07%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").Value = "'  -2.38%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").Value = "'  +3.57%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").Value = "'  +2.09%  "
$ws.Range("E51").Style = "Normal"
